$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("A8").Value = "SADRAEY"
$ws.Range("C8").Value = 6396.0
$ws.Range("D8").Value = -5.162918241106875

$ws.Range("A9").Value = "KROO"
$ws.Range("C9").Value = 6493.0
$ws.Range("D9").Value = -3.724644799797833

$ws.Range("A10").Value = "TORENBEEK_1976"
$ws.Range("C10").Value = 6240.0
$ws.Range("D10").Value = -7.4760177962018295

$ws.Range("A11").Value = "RAYMER"
$ws.Range("C11").Value = 6645.0
$ws.Range("D11").Value = -1.4708554897053137

$ws.Range("A12").Value = "TORENBEEK_2013"
$ws.Range("C12").Value = 7891.0
$ws.Range("D12").Value = 17.004285828553105

$ws.Range("A13").Value = "ROSKAM"
$ws.Range("C13").Value = 8148.0
$ws.Range("D13").Value = 20.814969069959535

$ws.Range("A14").Value = "JENKINSON"
$ws.Range("C14").Value = 9233.0
$ws.Range("D14").Value = 36.902873026869955

$ws.Range("A15").Value = "NICOLAI_1984"
$ws.Range("C15").Value = 6931.0
$ws.Range("D15").Value = 2.7698270279687693

$ws = $wb.Worksheets.Item("WING")
$ws.Range("A8").Value = "SADRAEY"
$ws.Range("C8").Value = 5801.0
$ws.Range("D8").Value = -18.85407281624436

$ws.Range("A9").Value = "KROO"
$ws.Range("C9").Value = 7124.0
$ws.Range("D9").Value = -0.3475977836450319

$ws.Range("A10").Value = "RAYMER"
$ws.Range("C10").Value = 8372.0
$ws.Range("D10").Value = 17.10975734914708

$ws.Range("A11").Value = "TORENBEEK_2013"
$ws.Range("C11").Value = 5858.0
$ws.Range("D11").Value = -18.056741692390876

$ws.Range("A12").Value = "TORENBEEK_1982"
$ws.Range("C12").Value = 6037.0
$ws.Range("D12").Value = -15.552842198184313

$ws.Range("A13").Value = "ROSKAM"
$ws.Range("C13").Value = 4280.0
$ws.Range("D13").Value = -40.13022438433475

$ws.Range("A14").Value = "JENKINSON"
$ws.Range("C14").Value = 886.0
$ws.Range("D14").Value = -87.60639691694406

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("A8").Value = "SADRAEY"
$ws.Range("C8").Value = 1040.0
$ws.Range("D8").Value = 34.092727831591574

$ws.Range("A9").Value = "NICOLAI_2013"
$ws.Range("C9").Value = 415.0
$ws.Range("D9").Value = -46.49184418258606

$ws.Range("A10").Value = "KROO"
$ws.Range("C10").Value = 738.0
$ws.Range("D10").Value = -4.845737365659056

$ws.Range("A11").Value = "RAYMER"
$ws.Range("C11").Value = 525.0
$ws.Range("D11").Value = -32.308959508090794

$ws.Range("A12").Value = "HOWE"
$ws.Range("C12").Value = 472.0
$ws.Range("D12").Value = -39.142531214893054

$ws.Range("A13").Value = "ROSKAM"
$ws.Range("C13").Value = 482.0
$ws.Range("D13").Value = -37.85317806266622

$ws.Range("A14").Value = "JENKINSON"
$ws.Range("C14").Value = 700.0
$ws.Range("D14").Value = -9.745279344121055

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("A8").Value = "SADRAEY"
$ws.Range("C8").Value = 765.0
$ws.Range("D8").Value = -1.3644838546465818

$ws.Range("A9").Value = "KROO"
$ws.Range("C9").Value = 497.0
$ws.Range("D9").Value = -35.91914833432595

$ws.Range("A10").Value = "RAYMER"
$ws.Range("C10").Value = 194.0
$ws.Range("D10").Value = -74.98654884679927

$ws.Range("A11").Value = "HOWE"
$ws.Range("C11").Value = 382.0
$ws.Range("D11").Value = -50.746709584934635

$ws.Range("A12").Value = "ROSKAM"
$ws.Range("C12").Value = 482.0
$ws.Range("D12").Value = -37.85317806266622

$ws.Range("A13").Value = "JENKINSON"
$ws.Range("C13").Value = 502.0
$ws.Range("D13").Value = -35.27447175821253

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("A10").Value = "KUNDU"
$ws.Range("C10").Value = 1389.0
$ws.Range("D10").Value = 116.79455344311013

$ws.Range("A11").Value = "ROSKAM"
$ws.Range("C11").Value = 1374.0
$ws.Range("D11").Value = 114.45335956143508

$ws.Range("A12").Value = "JENKINSON"
$ws.Range("C12").Value = 1410.0
$ws.Range("D12").Value = 120.07222487745521

$ws.Range("A17").Value = "KUNDU"
$ws.Range("C17").Value = 1389.0
$ws.Range("D17").Value = 116.79455344311013

$ws.Range("A18").Value = "ROSKAM"
$ws.Range("C18").Value = 1374.0
$ws.Range("D18").Value = 114.45335956143508

$ws.Range("A19").Value = "JENKINSON"
$ws.Range("C19").Value = 1410.0
$ws.Range("D19").Value = 120.07222487745521

$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("A11").Value = "TORENBEEK_1976"
$ws.Range("C11").Value = 2954.0
$ws.Range("D11").Value = 5.543652853730258

$ws.Range("A12").Value = "TORENBEEK_2013"
$ws.Range("C12").Value = 3458.0
$ws.Range("D12").Value = 23.551100733987553

$ws.Range("A13").Value = "KUNDU"
$ws.Range("C13").Value = 3265.0
$ws.Range("D13").Value = 16.655391525873153

$ws.Range("A18").Value = "TORENBEEK_1976"
$ws.Range("C18").Value = 2954.0
$ws.Range("D18").Value = 5.543652853730258

$ws.Range("A19").Value = "TORENBEEK_2013"
$ws.Range("C19").Value = 3458.0
$ws.Range("D19").Value = 23.551100733987553

$ws.Range("A20").Value = "KUNDU"
$ws.Range("C20").Value = 3265.0
$ws.Range("D20").Value = 16.655391525873153
